$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.624.72"
$ws.Range("E2").Value = "  -0.59%  "

$ws.Range("D3").Value = "2.506.06"
$ws.Range("E3").Value = "  +9.54%  "

$ws.Range("E4").Value = "  +0.80%  "

$ws.Range("D5").Formula = "'297.90"
$ws.Range("E5").Value = "  +0.28%  "

$ws.Range("D6").Formula = "'96.85"
$ws.Range("E6").Value = "  -2.09%  "

$ws.Range("D7").Formula = "'0.587"
$ws.Range("E7").Value = "  +3.40%  "

$ws.Range("E8").Value = "  +0.56%  "

$ws.Range("D9").Formula = "'0.541"
$ws.Range("E9").Value = "  +7.04%  "

$ws.Range("D10").Formula = "'36.75"
$ws.Range("E10").Value = "  +6.54%  "

$ws.Range("D11").Formula = "'0.0800"
$ws.Range("E11").Value = "  +3.04%  "

$ws.Range("D12").Formula = "'7.65"
$ws.Range("E12").Value = "  +9.52%  "

$ws.Range("D13").Value = "2.909.55"
$ws.Range("E13").Value = "  +10.78%  "

$ws.Range("E14").Value = "  +2.46%  "

$ws.Range("D15").Value = "2.529.93"
$ws.Range("E15").Value = "  +9.64%  "

$ws.Range("D16").Formula = "'0.873"
$ws.Range("E16").Value = "  +10.31%  "

$ws.Range("D17").Formula = "'14.48"
$ws.Range("E17").Value = "  +7.20%  "

$ws.Range("D18").Value = "45.784.23"
$ws.Range("E18").Value = "  -0.05%  "

$ws.Range("D19").Formula = "'13.19"
$ws.Range("E19").Value = "  +6.20%  "

$ws.Range("D20").Value = "0.0₃0968"
$ws.Range("E20").Value = "  +1.49%  "

$ws.Range("E21").Value = "  +12.21%  "

$ws.Range("D22").Formula = "'68.57"
$ws.Range("E22").Value = "  +5.27%  "

$ws.Range("D23").Formula = "'249.64"
$ws.Range("E23").Value = "  +2.75%  "

$ws.Range("D24").Formula = "'2.86"
$ws.Range("E24").Value = "  +2.94%  "

$ws.Range("D25").Formula = "'2.04"
$ws.Range("E25").Value = "  +9.84%  "

$ws.Range("D26").Formula = "'0.999"
$ws.Range("E26").Value = "  -0.35%  "

$ws.Range("D27").Formula = "'40.34"
$ws.Range("E27").Value = "  +0.35%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Formula = "'22.90"
$ws.Range("E28").Value = "  +14.57%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Formula = "'10.04"
$ws.Range("E29").Value = "  +5.62%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Formula = "'2.24"
$ws.Range("E30").Value = "  +1.43%  "

$ws.Range("D31").Formula = "'3.79"
$ws.Range("E31").Value = "  +12.51%  "

$ws.Range("D32").Formula = "'5.74"
$ws.Range("E32").Value = "  +8.51%  "

$ws.Range("D33").Formula = "'2.78"
$ws.Range("E33").Value = "  -0.55%  "

$ws.Range("D34").Formula = "'2.16"
$ws.Range("E34").Value = "  +30.11%  "

$ws.Range("D35").Formula = "'148.86"
$ws.Range("E35").Value = "  +3.73%  "

$ws.Range("D36").Formula = "'0.0803"
$ws.Range("E36").Value = "  +5.01%  "

$ws.Range("E37").Value = "  +5.09%  "

$ws.Range("E38").Value = "  +3.06%  "

$ws.Range("D39").Formula = "'15.58"
$ws.Range("E39").Value = "  +1.72%  "

$ws.Range("D40").Formula = "'4.08"
$ws.Range("E40").Value = "  +7.50%  "

$ws.Range("D41").Formula = "'0.0310"
$ws.Range("E41").Value = "  +5.70%  "

$ws.Range("D42").Formula = "'3.36"
$ws.Range("E42").Value = "  +8.85%  "

$ws.Range("D43").Value = "2.001.93"
$ws.Range("E43").Value = "  +10.72%  "

$ws.Range("D44").Formula = "'1.00"
$ws.Range("E44").Value = "  +0.28%  "

$ws.Range("D45").Formula = "'90.70"
$ws.Range("E45").Value = "  -2.93%  "

$ws.Range("D46").Formula = "'16.59"
$ws.Range("E46").Value = "  +34.45%  "

$ws.Range("D47").Formula = "'1.78"
$ws.Range("E47").Value = "  -3.38%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Formula = "'8.79"
$ws.Range("E48").Value = "  +13.47%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Formula = "'105.38"
$ws.Range("E49").Value = "  +11.94%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Formula = "'0.194"
$ws.Range("E50").Value = "  +5.96%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.762.80"
$ws.Range("E51").Value = "  +10.21%  "
